$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 141 - this shifts the existing rows
# 141-158 down to 143-160, preserving all of their data untouched.
$ws.Rows("141:142").Insert()

# Row 141 - new weekly entry (Primera), matching the next data point
# for Betarraga (4-unit package).
$ws.Range("A141").Value = 1
$ws.Range("B141").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C141").Value = "Arica y Parinacota"
$ws.Range("D141").Value = 44449
$ws.Range("E141").Value = 15
$ws.Range("F141").Value = 100114014
$ws.Range("G141").Value = "Betarraga"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 1200
$ws.Range("K141").Value = 500
$ws.Range("L141").Value = 550
$ws.Range("M141").Value = 525
$ws.Range("N141").Value = "$/paquete 4 unidades"
$ws.Range("O141").Value = "Región de Arica y Parinacota"
$ws.Range("P141").Value = 131
$ws.Range("Q141").Value = 4
$ws.Range("R141").Value = "Hortaliza"

# Row 142 - same week (Segunda), 5-unit package.
$ws.Range("A142").Value = 1
$ws.Range("B142").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C142").Value = "Arica y Parinacota"
$ws.Range("D142").Value = 44449
$ws.Range("E142").Value = 15
$ws.Range("F142").Value = 100114014
$ws.Range("G142").Value = "Betarraga"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Segunda"
$ws.Range("J142").Value = 1200
$ws.Range("K142").Value = 500
$ws.Range("L142").Value = 550
$ws.Range("M142").Value = 525
$ws.Range("N142").Value = "$/paquete 5 unidades"
$ws.Range("O142").Value = "Región de Arica y Parinacota"
$ws.Range("P142").Value = 105
$ws.Range("Q142").Value = 5
$ws.Range("R142").Value = "Hortaliza"
